$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Date" header added in A1 (previously empty)
$ws.Range("A1").Value = "Date"

# A2's date stamp loses its leading underscore, but must remain text
# (not be reinterpreted as a number) - force with a leading apostrophe,
# same as typing it directly into Excel.
$ws.Range("A2").Value = "'20240226"

# Columns C:G are reordered (site codes rearranged); keep each header's
# value travelling with it to its new column.
$ws.Range("C1").Value = "PKV"
$ws.Range("D1").Value = "STL"
$ws.Range("E1").Value = "THL"
$ws.Range("F1").Value = "AZC"
$ws.Range("G1").Value = "SDU"

$ws.Range("C2").Value = 159
$ws.Range("D2").Value = 68
$ws.Range("E2").Value = 110
$ws.Range("F2").Value = 220
$ws.Range("G2").Value = 39
